$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new year column S (2022) mirroring the formatting of column R
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S5").Value = 13.5

# Updated existing values in row 5
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.9

# Update the active selection to match the new extent
$ws.Range("S7:S8").Select()
